# commit for loginPage module
#
# Inserts a new "TC001_LoginPageValidation" row right after the header
# row (pushing the three existing test cases down by one row and
# renumbering them TC002.."TC004"), then repoints the e-mail/password
# hyperlinks at their (shifted) cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Shift the three existing data rows down one (copies both the
#        value AND the cell formatting/style), working bottom-up so a
#        row is never overwritten before it has been copied forward. ---
$ws.Range("A4").Copy($ws.Range("A5"))
$ws.Range("B4").Copy($ws.Range("B5"))
$ws.Range("C4").Copy($ws.Range("C5"))

$ws.Range("A3").Copy($ws.Range("A4"))
$ws.Range("B3").Copy($ws.Range("B4"))
$ws.Range("C3").Copy($ws.Range("C4"))

$ws.Range("A2").Copy($ws.Range("A3"))
$ws.Range("B2").Copy($ws.Range("B3"))
$ws.Range("C2").Copy($ws.Range("C3"))

# --- 2. New row 2: login-page validation test case, email/password blank ---
$ws.Range("A2").Value = "TC001_LoginPageValidation"
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""

# --- 3. Renumber the shifted test cases ---
$ws.Range("A3").Value = "TC002_invalidemail"
$ws.Range("A4").Value = "TC003_invalidpassword"
$ws.Range("A5").Value = "TC004_validlogin"

# --- 4. Rebuild the hyperlinks so each e-mail/password cell points at the
#        right mailto: target again (row 2 now has no hyperlinks). ---
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:gowtham@yahoo.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:test@123")
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:abc@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:Test@123")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:lotica.aitech@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:Lotica@123")

# --- 5. Match the author's final selection ---
$ws.Range("B11").Select()
